$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = 'magapoke_2026-02-11'

$ws.Range("A1").Value = 'rank'
$ws.Range("B1").Value = 'title'
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'ハンドレッドノート－アグリーダック－'
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！'
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'ペンの夢に紅をさす'
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = '春くらり'
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = '黒月のイェルクナハト'
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = '黄昏町プリズナーズ'
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'お前がヤったんだろ！'
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = '【爆アド】生まれた直後から最強悪霊と脳内バトルしてたら魔力量が測定可能域を超えてました〜悪憑の子の謙虚な覇道〜'
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'K-9~警視庁公安部公安第9課異能対策係~'
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 'アイドラトリィ'
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'あの島の海音荘'
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'せいぶつ部の田辺くん'
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'ハードワーカー中田'
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'ゼロとヒャク'
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 'ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜'
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = '邪目さんは邪神です'
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 'その青春'
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 'ドリーム☆ジャンボ☆ガール'
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = '屋根の下のアルテミス'
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = '普通の本はありません！'
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'MYS'
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = '明智ナンバーワン'
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = '白鳥運子は31画'
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = '歪みの虜'
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 'わが投資術　市場は誰に微笑むか'
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = '限界集落を脱村した錬金術士、都会で"最強"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～'
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = '人生逆転ダンジョン'
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = '異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～'
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = '宇曽田みのりの代用料理'
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = '篝家の８兄弟'
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = '君が監督！'
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'ともだちづくり'
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = '皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～'
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 'ナキナギ'
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = '鳴るさんだぁ'
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'ナマイキ旭ちゃんをわからせたい'
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'じゅーくぼっくす'
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = '眠れる森のレガ'
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = '追放されなかった男　～二度目の人生は土下座から始まりました～'
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = '平成転生'
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = '白銀のキュイジーヌ～明治外交官の料理人～'
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = 'GURU'
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = '永久のユウグレ'
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 'ハプスブルク家の華麗なる受難'
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 'ch登録お願いします！'
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = '夜鐘のキト'
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = 'イエティ、とある日々'
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = 'JK Biker'
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = '〈小市民〉 春期限定いちごタルト事件'
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = '花子狩り'
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = 'きゃわるり方程式'

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Move($null, $last)
